# "add comments to functions" -- the underlying xlsx re-save mostly carries
# view/formatting deltas (column widths, selection, scroll position) rather
# than cell-comment content. Re-apply the reachable parts of that re-save
# through the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------------
# Column A keeps its original best-fit width, so it is intentionally left
# untouched (re-running AutoFit on it would strip the best-fit flag).
# Columns B, C, D, E, F and H were manually widened (no longer "best fit"),
# so set explicit widths for them. ColumnWidth is specified in characters.
$ws.Columns.Item(2).ColumnWidth = 14.333333333333334   # B -> ~15.14 chars
$ws.Columns.Item(3).ColumnWidth = 79.66666666666667    # C -> ~80.43 chars
$ws.Columns.Item(4).ColumnWidth = 57.5                 # D -> ~58.29 chars
$ws.Columns.Item(5).ColumnWidth = 53.666666666666664   # E -> ~54.57 chars
$ws.Columns.Item(6).ColumnWidth = 48.666666666666664   # F -> ~49.43 chars
$ws.Columns.Item(8).ColumnWidth = 56.666666666666664   # H -> ~57.57 chars

# --- View state: scroll position + selection -----------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1

$ws.Range("H12").Select()
